$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("BC4").Value = 126
$ws.Range("BD4").Value = 126

# Row 6
$ws.Range("G6").Value = 1.9
$ws.Range("I6").Value = 3.6
$ws.Range("J6").Value = 2.4
$ws.Range("L6").Value = 4
$ws.Range("X6").Value = 11
$ws.Range("AC6").Value = 17
$ws.Range("AK6").Value = 26
$ws.Range("AO6").Value = 9.5
$ws.Range("AS6").Value = 81
$ws.Range("AX6").Value = 19

# Row 9
$ws.Range("G9").Value = 1.55
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("Q9").Value = 1.95
$ws.Range("R9").Value = 1.95
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 2.75
$ws.Range("AB9").Value = 29
$ws.Range("AE9").Value = 17
$ws.Range("AO9").Value = 8
$ws.Range("AP9").Value = 21
$ws.Range("AT9").Value = 2.75

# Row 12
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 13

# Row 19
$ws.Range("I19").Value = 4.5
$ws.Range("M19").Value = 1.14
$ws.Range("N19").Value = 5.5
$ws.Range("Q19").Value = 2.88
$ws.Range("R19").Value = 1.4
$ws.Range("W19").Value = 5
$ws.Range("AC19").Value = 5.5
$ws.Range("AI19").Value = 17
$ws.Range("AK19").Value = 41
$ws.Range("AL19").Value = 51
$ws.Range("AP19").Value = 29
$ws.Range("AU19").Value = 10

# Row 20
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 7.4
$ws.Range("I20").Value = 1.08
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 3.45
$ws.Range("L20").Value = 1.32
$ws.Range("Q20").Value = 1.26
$ws.Range("T20").Value = 4.7
$ws.Range("U20").Value = 2.1
$ws.Range("V20").Value = 1.65
$ws.Range("W20").Value = 120
$ws.Range("X20").Value = 500
$ws.Range("Y20").Value = 100
$ws.Range("AA20").Value = 600
$ws.Range("AB20").Value = 250
$ws.Range("AC20").Value = 26
$ws.Range("AD20").Value = 21
$ws.Range("AF20").Value = 150
$ws.Range("AH20").Value = 7.6
$ws.Range("AI20").Value = 12.5
$ws.Range("AJ20").Value = 6.8
$ws.Range("AK20").Value = 10.75
$ws.Range("AL20").Value = 32
$ws.Range("AN20").Value = 24
$ws.Range("AO20").Value = 175
$ws.Range("AP20").Value = 80
$ws.Range("AT20").Value = 4.7
$ws.Range("AU20").Value = 10.5
$ws.Range("AV20").Value = 70
$ws.Range("AX20").Value = 4.2
$ws.Range("AZ20").Value = 7.3
